$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column B holds numeric-looking codes that must remain stored as text
# (matches source data format), so force text number format first.
$ws.Range("B2:B10").NumberFormat = "@"

# Row 2
$ws.Range("A2").Value = 'favoravel'
$ws.Range("B2").Value = '11'
$ws.Range("C2").Value = '"Racismo reverso não existe"'
$ws.Range("D2").Value = 175
$ws.Range("E2").Value = 363519
$ws.Range("F2").Value = 53344
$ws.Range("G2").Value = 6503
$ws.Range("H2").Value = 3851
$ws.Range("I2").Value = 'BabuSantana, brunogagliasso, YuriMarcal, brunoformiga, essediafoilouco, iglesbiteriana, thiamparo, HistoriaNoPaint, rodrigocapelo, MussumAlive, cinefilo_K, andrefran, LeviKaique, dominiopop, fabiunascimento, Ticostacruz, taliriapetrone, SimoneEhNois, QuebrandoOTabu, startupdareal, aanonnyma, AfroSailor, lolaescreva, ma_azevedo94, Cauefabiano, joelluiz_adv, franca_rodrigo, ReTintaPreta, danielleonawale, knjcheeks_, gaby17pontes, GabiRAR1, rogercipo, davidmirandario, lgqueiroga, heroinadolixo, gleisi, Savagefiction, Pessoa_Cansada, carolizaaando, Rosy_Oliveiira, danvieirass, jpgadelhaof, jojopancada, soutamires_sp, soutoverso, pestodeboldo, LeonelRadde, levieiraprof, MidiaNINJA, caduadv, FellerMarcelo, MarceloUchoa_, EstreLaaany, TeresaCristina, michelpinho, tiamaoficial, luyarafranco, patriers, iamKalera, jhonpaim15, afrocrente, Nailahnv, sigaoflavio, sunamita_nobre, jaoseupimenta, hospicio_brasil, TigerGames, HeslaineVieira, livialaranjeira, Crissy_98_, profsosa13, religiaosincera, marioadolfo, esquerdeando, gdoweber, joaofelipenobre, eueoyoon, _rosaneborges, anunesrocha, MariluPamc, buerolol, souljazzca, realodara, livrosdodrii, jfmathias, afro_hey, vanessasoaresc4, mayrasigwalt, waltermoraes_, lucas_kurz, jaciarabri, LianaCirne, Cris_Barbieri, passarosErosas, brauneoficial, blckklucas, donairene13, geekcom2, oBrunoRomano, Sucubus, refugefefo, SirLucasMatheus, eurickrodrigues, evandrof, AeroportoD, florapaulita, danibacedo, Daniel07091992, jonasdiandrade, beatxriz, jgprates, afroestima2, ChicaoBulhoes, allisonaw1965, orlandoguerreir, rafhaelnep, observint, JeGiacometInda, OManoRogerio, joaooribeiro26, brunnosarttori, drikbarbosa, Pablo_Peixoto, takemeout, monica_benicio, eenrietti, BrunoCostoli, mota97fm, erahsfeliz, desconstrunutri, ovictorjame, gobletofpjo, eusamantalima, bea_brazx, slc_cavalcante, adalto_edno, Sybylla_, mdmdaiane, 1cesgusto, souarthurlima, FredLAFernandes, guirocha82, ashleymlia, llcncl, bragacamila1'

# Row 3
$ws.Range("A3").Value = 'favoravel'
$ws.Range("B3").Value = '12'
$ws.Range("C3").Value = '"O trainee magalu é uma ação afirmativa desejável"'
$ws.Range("D3").Value = 110
$ws.Range("E3").Value = 94447
$ws.Range("F3").Value = 10649
$ws.Range("G3").Value = 2959
$ws.Range("H3").Value = 1431
$ws.Range("I3").Value = 'thiamparo, danielagomesphd, slpng_giants_pt, DTangerinoPenal, cynaramenezes, jeff_theblack, LeviKaique, QuebrandoOTabu, AndreaMPacha, Savagefiction, tesoureiros, hudsonbonatto, mariliadf2, dasilvabenedita, RMafei, Neka_BR, Ticostacruz, brunnosarttori, Omardeideais, AdrianaCarranca, luisaciteli, GabrielCBrasil, davicalbez, RafaelValim7, NA, fernandapsol, bolsoregrets, srlm, tiamaoficial, luziel__, Sen_Cristovam, andrefran, andrefatala, chambzrs, edufelipe101, agrdeumsm, SamPancher, heelobrandao, marcogomes, anacarla_abrao, franca_rodrigo, passalanorh, GabineteOcio, pbdjulia, TarcisoRenova, potenciasnegras, TatiNefertari, oviniporto, desisalg, isantanax, wendy_andrade, barretonessa, albertocalmeida, andersonsenaxxp, jennieyees, AliQuintiliano, emirsader, jojopancada, rogercipo, Bebeto_Esposito, felip0c, JaumGodoy, Sonia_Cout_, Joao_Gini, lapena, Iberedias, dz7king, mailsonmcj, gamerpobre12344, FamosoLucas_, jgprates, AndreGomesF, luclsluiz, Legurbano, joelluiz_adv, antonionetopdt, EuSouLivres, buruaca, RicardoWeber, teeaggo, naty_andradde, gabrielzep, juliamolusco7, omgerva, kauancoellho, medoedeliriobr, tainadepaularj, SawaraKali, Flavio_Sampaio, _nelsoncezar, Tiagoonie, cruz_elianalves, nathali20044244, ago_almeida, profaflavia, talitismo, hilde_angel'

# Row 4
$ws.Range("A4").Value = 'contrario'
$ws.Range("B4").Value = '14'
$ws.Range("C4").Value = '"O trainee magalu promove o racismo"'
$ws.Range("D4").Value = 71
$ws.Range("E4").Value = 80522
$ws.Range("F4").Value = 14414
$ws.Range("G4").Value = 5736
$ws.Range("H4").Value = 2012
$ws.Range("I4").Value = 'carlosjordy, LorenzonItalo, profpaulamarisa, LuizCamargoVlog, spinellirio, depheliolopes, josuenunes, BlackDogBC, Desesquerdizada, Jouberth19, romollerSP, rmotta2, FernandoHoliday, Bolsoneas, RRDECA_, opropriofaka, dezacrvg, CrisMenshova, AlessandroLoio2, simmer_lara, CanaldoNegaoo, LucianaSV29, PastorLiomar, _charizard100, canalCCore2, WagnerThomazoni, thaispsic, JulioOliSantos2, ericayhwh, Sirlene_Emanuel, jaohff, ribas1960, PPReacaFla2, GauchaLih, NA, fabiomello1010, CaioCarlosibg, erickirios, guilhermedecnop, paulocruzphi, teresinhalopes, gabrielferna_a, Arthurdinizrd, kkgbraga, jesus_filha, VlogdoLisboa, FernandoMessina, doxxxx, viniciussexto, ToniTosti1, _VF20, FredRC, MonicaMachado38, nerdclassico, WolfConservador, ajulysantos, rinaldidigilio, thiagosiqueira5'

# Row 5
$ws.Range("A5").Value = 'outros'
$ws.Range("B5").Value = '98'
$ws.Range("C5").Value = '"Matérias jornalísticas"'
$ws.Range("D5").Value = 53
$ws.Range("E5").Value = 39630
$ws.Range("F5").Value = 4327
$ws.Range("G5").Value = 4407
$ws.Range("H5").Value = 1601
$ws.Range("I5").Value = 'folha, exame, JornalOGlobo, conexaopolitica, revistaforum, jornalextra, brasil247, JornalDaCidadeO, Estadao, BlogdoNoblat, UOLNoticias, revistaoeste, congressoemfoco, DiarioPE, UOL, DCM_online, Metropoles, gazetadopovo, CNNBrBusiness, bbcbrasil, elpais_brasil, jc_pe, madeleinelacsko, correio, laurojardim, em_com, flaviaol'

# Row 6
$ws.Range("A6").Value = 'outros'
$ws.Range("B6").Value = '99'
$ws.Range("C6").Value = '"Posicionamento não identificado"'
$ws.Range("D6").Value = 38
$ws.Range("E6").Value = 11416
$ws.Range("F6").Value = 840
$ws.Range("G6").Value = 568
$ws.Range("H6").Value = 498
$ws.Range("I6").Value = 'andrezadelgado, jpgadelhaof, canaldasbee, luide, MussumAlive, paulocruzphi, ittsquel, alisc, GirassolRafa, nathali20044244, FredRC, dedehcamargo, Savagefiction, tesoureiros, brunanarcizo, femisapien_z, ChampMargareth, luccaoneal, TarcisoRenova, JuanSavedra_, souarthurlima, _giovanirocha, artedaguerracnl, IzaVicent, sourodrii, dannielduque, JacyCarvalho, Rayctjay, Lethiscya, reclamiranda, wasabinoolho, joaolordelo, rdfmedeiros, afroestima2, biraiorio'

# Row 7
$ws.Range("A7").Value = 'contrario'
$ws.Range("B7").Value = '13'
$ws.Range("C7").Value = '"O trainee magalu desrespeita a meritocracia"'
$ws.Range("D7").Value = 12
$ws.Range("E7").Value = 8644
$ws.Range("F7").Value = 1567
$ws.Range("G7").Value = 1932
$ws.Range("H7").Value = 433
$ws.Range("I7").Value = 'FernandoHoliday, phillipGlotok, ClaudeLuca_, vigilantshitter, opropriofaka, edibertoalves, junggukkah, BettoFerreira6, frz_daniel, Juh_Oliveira_I, luisfred63'

# Row 8
$ws.Range("A8").Value = 'contrario'
$ws.Range("B8").Value = '17'
$ws.Range("C8").Value = '"Posicionamentos contrários variados"'
$ws.Range("D8").Value = 16
$ws.Range("E8").Value = 8008
$ws.Range("F8").Value = 998
$ws.Range("G8").Value = 510
$ws.Range("H8").Value = 162
$ws.Range("I8").Value = 'hoc111, paulocruzphi, rmotta2, Felippe_Hermes, CrysthianeA, thaispsic, mattlimn, gramich, canalCCore2, m_blazar, _Renato, EKrominski, TchiMad, de_botequim'

# Row 9
$ws.Range("A9").Value = 'contrario'
$ws.Range("B9").Value = '15'
$ws.Range("C9").Value = '"O trainee magalu realiza um ''apartheid'' na sociedade"'
$ws.Range("D9").Value = 11
$ws.Range("E9").Value = 3706
$ws.Range("F9").Value = 596
$ws.Range("G9").Value = 291
$ws.Range("H9").Value = 123
$ws.Range("I9").Value = 'nerdclassico, carlinhoscury, LucianaSV29, SantoElayne2020, Nivea_SpesEst, carlaguimaraes7, CRCF_84, MarizMarcella, mjmacul_lima, paulocnf, fernando_g_f_'

# Row 10
$ws.Range("A10").Value = 'contrario'
$ws.Range("B10").Value = '16'
$ws.Range("C10").Value = '"Aceitação apenas de critérios universais"'
$ws.Range("D10").Value = 10
$ws.Range("E10").Value = 3504
$ws.Range("F10").Value = 476
$ws.Range("G10").Value = 152
$ws.Range("H10").Value = 39
$ws.Range("I10").Value = 'nerdclassico, lpdossj, LucianaSV29, letparks, sampaio19091, julioreis78, Alpargatas11, Buzz91272883, moneymakerbr'
